# "Generate Report for Handback" — refresh the handback timestamps / status
# recorded by the localization report generator.
#
# The workbook has three tables:
#   "Overview" : one row per source file, "Latest HO Xliff Generate Date"
#   "zh-cn"    : per-locale detail for the zh-cn handback
#   "de-de"    : per-locale detail for the de-de handback
#
# Both data rows (2 and 3) of each sheet share the same timestamps / status
# since they were generated by the same report run, so every value below is
# written to both rows.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-14 07:29:34"
$wsOverview.Range("G3").Value = "2016-11-14 07:29:34"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority (column E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# Correspond Handoff Datetime (column H)
$wsZhCn.Range("H2").Value = "2016-11-14 07:29:21"
$wsZhCn.Range("H3").Value = "2016-11-14 07:29:21"
# Correspond Handback DateTime (column K)
$wsZhCn.Range("K2").Value = "2016-11-14 07:30:16"
$wsZhCn.Range("K3").Value = "2016-11-14 07:30:16"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority (column E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# Correspond Handoff Datetime (column H)
$wsDeDe.Range("H2").Value = "2016-11-14 07:29:34"
$wsDeDe.Range("H3").Value = "2016-11-14 07:29:34"
# Correspond Handback DateTime (column K)
$wsDeDe.Range("K2").Value = "2016-11-14 07:30:34"
$wsDeDe.Range("K3").Value = "2016-11-14 07:30:34"
